$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 428.2
$ws.Range("I4").Value = 428.2
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 428.2
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -314.2
$ws.Range("N4").Value = ""
$ws.Range("H87").Value = 42609
$ws.Range("J87").Value = 42609
$ws.Range("L87").Value = 42609
$ws.Range("N87").Value = -45105
$ws.Range("H90").Value = 42609
$ws.Range("J90").Value = 42609
$ws.Range("L90").Value = 127827
$ws.Range("N90").Value = -140307
$ws.Range("H113").Value = 3389.4443
$ws.Range("I113").Value = 2712.25
$ws.Range("J113").Value = 3931.2
$ws.Range("K113").Value = 2712.25
$ws.Range("L113").Value = 3931.2
$ws.Range("M113").Value = 541.75
$ws.Range("N113").Value = -10439.2

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3976.7812
$ws.Range("I32").Value = 3476.4238
$ws.Range("K32").Value = 3476.4238
$ws.Range("M32").Value = -3189.4238
$ws.Range("H37").Value = 29333.334
$ws.Range("J37").Value = 29333.334
$ws.Range("L37").Value = 29333.334
$ws.Range("N37").Value = -29879.334
$ws.Range("H44").Value = 89500
$ws.Range("J44").Value = 89500
$ws.Range("L44").Value = 89500
$ws.Range("N44").Value = -90476
$ws.Range("H61").Value = 7331.32
$ws.Range("I61").Value = 7603.609
$ws.Range("J61").Value = 4200
$ws.Range("K61").Value = 7603.609
$ws.Range("L61").Value = 4200
$ws.Range("M61").Value = -7391.609
$ws.Range("N61").Value = -4624
$ws.Range("H80").Value = 48799.2
$ws.Range("J80").Value = 48799.2
$ws.Range("L80").Value = 48799.2
$ws.Range("N80").Value = -50795.2
$ws.Range("H83").Value = 48799.2
$ws.Range("J83").Value = 48799.2
$ws.Range("L83").Value = 146397.6
$ws.Range("N83").Value = -156381.6
$ws.Range("H110").Value = 1735.9286
$ws.Range("I110").Value = 1791.9166
$ws.Range("J110").Value = 1400
$ws.Range("K110").Value = 1791.9166
$ws.Range("L110").Value = 1400
$ws.Range("M110").Value = 253.0834
$ws.Range("N110").Value = -5490
$ws.Range("H122").Value = 15026.24
$ws.Range("I122").Value = 1896.5883
$ws.Range("J122").Value = 42926.75
$ws.Range("K122").Value = 5689.7649
$ws.Range("L122").Value = 128780.25
$ws.Range("M122").Value = -3239.7649
$ws.Range("N122").Value = -133680.25
$ws.Range("H136").Value = 7331.32
$ws.Range("I136").Value = 7603.609
$ws.Range("J136").Value = 4200
$ws.Range("K136").Value = 22810.827
$ws.Range("L136").Value = 12600
$ws.Range("M136").Value = -20260.827
$ws.Range("N136").Value = -17700
$ws.Range("H139").Value = 114983.875
$ws.Range("J139").Value = 114983.875
$ws.Range("L139").Value = 114983.875
$ws.Range("N139").Value = -125263.875

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 2359.6
$ws.Range("I54").Value = 1449.75
$ws.Range("J54").Value = 5999
$ws.Range("K54").Value = 1449.75
$ws.Range("L54").Value = 5999
$ws.Range("M54").Value = -965.75
$ws.Range("N54").Value = -6967
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = ""
$ws.Range("H82").Value = 22766.166
$ws.Range("I82").Value = 9400
$ws.Range("K82").Value = 9400
$ws.Range("M82").Value = -9017
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = ""
$ws.Range("H85").Value = 22766.166
$ws.Range("I85").Value = 9400
$ws.Range("K85").Value = 9400
$ws.Range("M85").Value = -8074
$ws.Range("H107").Value = 3024.5
$ws.Range("I107").Value = 2138.6
$ws.Range("J107").Value = 3910.4
$ws.Range("K107").Value = 2138.6
$ws.Range("L107").Value = 3910.4
$ws.Range("M107").Value = -218.5999999999999
$ws.Range("N107").Value = -7750.4
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1616.5454
$ws.Range("I31").Value = 1616.5454
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1616.5454
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1321.5454
$ws.Range("N31").Value = ""
$ws.Range("H34").Value = 1616.5454
$ws.Range("I34").Value = 1616.5454
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1616.5454
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1414.5454
$ws.Range("N34").Value = ""
$ws.Range("H50").Value = 51666.668
$ws.Range("J50").Value = 49000
$ws.Range("L50").Value = 49000
$ws.Range("N50").Value = -50250
$ws.Range("H51").Value = 80000
$ws.Range("J51").Value = 80000
$ws.Range("L51").Value = 80000
$ws.Range("N51").Value = -81472
$ws.Range("H58").Value = 2088.8635
$ws.Range("I58").Value = 1945.579
$ws.Range("J58").Value = 2996.3333
$ws.Range("K58").Value = 1945.579
$ws.Range("L58").Value = 2996.3333
$ws.Range("M58").Value = -1742.579
$ws.Range("N58").Value = -3402.3333
$ws.Range("H60").Value = 46618.6
$ws.Range("J60").Value = 61000
$ws.Range("L60").Value = 61000
$ws.Range("N60").Value = -62022
$ws.Range("H61").Value = 80000
$ws.Range("J61").Value = 80000
$ws.Range("L61").Value = 80000
$ws.Range("N61").Value = -80696
$ws.Range("H132").Value = 8978.457
$ws.Range("I132").Value = 4912.185
$ws.Range("J132").Value = 22702.125
$ws.Range("K132").Value = 14736.555
$ws.Range("L132").Value = 68106.375
$ws.Range("M132").Value = -12206.555
$ws.Range("N132").Value = -73166.375
$ws.Range("H136").Value = 2088.8635
$ws.Range("I136").Value = 1945.579
$ws.Range("J136").Value = 2996.3333
$ws.Range("K136").Value = 5836.737
$ws.Range("L136").Value = 8988.999899999999
$ws.Range("M136").Value = -3286.737
$ws.Range("N136").Value = -14088.9999

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 7500
$ws.Range("J31").Value = 7500
$ws.Range("L31").Value = 22500
$ws.Range("N31").Value = -23076
$ws.Range("H114").Value = 1214.1666
$ws.Range("I114").Value = 1257.4
$ws.Range("J114").Value = 998
$ws.Range("K114").Value = 3772.2
$ws.Range("L114").Value = 2994
$ws.Range("M114").Value = -518.2000000000003
$ws.Range("N114").Value = -9502
$ws.Range("H137").Value = 10276.235
$ws.Range("I137").Value = 6870.1665
$ws.Range("J137").Value = 12134.091
$ws.Range("K137").Value = 20610.4995
$ws.Range("L137").Value = 36402.273
$ws.Range("M137").Value = -15510.4995
$ws.Range("N137").Value = -46602.273
$ws.Range("H139").Value = 3801.842
$ws.Range("I139").Value = 1698.2142
$ws.Range("K139").Value = 5094.642599999999
$ws.Range("M139").Value = 45.35740000000078
$ws.Range("H141").Value = 15058.375
$ws.Range("I141").Value = 15058.375
$ws.Range("K141").Value = 45175.125
$ws.Range("M141").Value = -39995.125

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4071
$ws.Range("I102").Value = 3574.4075
$ws.Range("K102").Value = 3574.4075
$ws.Range("M102").Value = -1952.4075
$ws.Range("H107").Value = 382
$ws.Range("I107").Value = 529.2222
$ws.Range("J107").Value = 161.16667
$ws.Range("K107").Value = 529.2222
$ws.Range("L107").Value = 161.16667
$ws.Range("M107").Value = 1390.7778
$ws.Range("N107").Value = -4001.16667
$ws.Range("H126").Value = 4550.2354
$ws.Range("I126").Value = 4248.724
$ws.Range("J126").Value = 6299
$ws.Range("K126").Value = 12746.172
$ws.Range("L126").Value = 18897
$ws.Range("M126").Value = -10276.172
$ws.Range("N126").Value = -23837
$ws.Range("H132").Value = 6198.5757
$ws.Range("I132").Value = 4744.635
$ws.Range("J132").Value = 11598.929
$ws.Range("K132").Value = 14233.905
$ws.Range("L132").Value = 34796.787
$ws.Range("M132").Value = -11703.905
$ws.Range("N132").Value = -39856.787

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4500.5713
$ws.Range("I40").Value = 4004.7083
$ws.Range("K40").Value = 4004.7083
$ws.Range("M40").Value = -3868.7083
$ws.Range("H132").Value = 42819.55
$ws.Range("I132").Value = 45976.137
$ws.Range("J132").Value = 3888.3333
$ws.Range("K132").Value = 137928.411
$ws.Range("L132").Value = 11664.9999
$ws.Range("M132").Value = -135398.411
$ws.Range("N132").Value = -16724.9999
